$p = $ppt.ActivePresentation

# Slide 1: title "First" + " " + "slide" -> "First " + "slide"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 6).Text = "First "

# Slide 3: title "Third" + " " + "slide" -> "Third " + "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 6).Text = "Third "
